$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'291.20"
$ws.Range("E2").Value = "'-1.32%"
$ws.Range("D3").Value = "'30.81"
$ws.Range("E3").Value = "'-1.11%"
$ws.Range("D4").Value = "'4.901"
$ws.Range("E4").Value = "'-1.10%"
$ws.Range("D5").Value = "'0.07248"
$ws.Range("E5").Value = "'-1.51%"
$ws.Range("D6").Value = "'2.328"
$ws.Range("E6").Value = "'26.33%"
$ws.Range("D7").Value = "'7.676"
$ws.Range("E7").Value = "'-0.06%"
$ws.Range("D8").Value = "'3.707"
$ws.Range("E8").Value = "'-1.18%"
$ws.Range("D9").Value = "'0.8964"
$ws.Range("E9").Value = "'-1.60%"
$ws.Range("D10").Value = "'0.1666"
$ws.Range("E10").Value = "'-0.21%"
$ws.Range("D11").Value = "'0.07919"
$ws.Range("E11").Value = "'3.69%"
$ws.Range("D12").Value = "'0.08015"
$ws.Range("E12").Value = "'-1.37%"
$ws.Range("D13").Value = "'0.03107"
$ws.Range("E13").Value = "'4.24%"
$ws.Range("D14").Value = "'0.1003"
$ws.Range("E14").Value = "'0.50%"
$ws.Range("D15").Value = "'0.001497"
$ws.Range("E15").Value = "'0.23%"
$ws.Range("D16").Value = "'0.005863"
$ws.Range("E16").Value = "'4.03%"
$ws.Range("D17").Value = "'3.465"
$ws.Range("E17").Value = "'0.15%"
$ws.Range("E18").Value = "'-1.13%"
$ws.Range("D19").Value = "'0.3319"
$ws.Range("E19").Value = "'1.43%"
$ws.Range("E20").Value = "'-0.54%"
$ws.Range("D21").Value = "'3.966"
$ws.Range("E21").Value = "'-8.55%"
$ws.Range("D22").Value = "'0.2202"
$ws.Range("E22").Value = "'10.29%"
$ws.Range("D23").Value = "'0.04530"
$ws.Range("E23").Value = "'1.27%"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("E24").Value = "'-1.03%"
$ws.Range("D25").Value = "'0.004412"
$ws.Range("E25").Value = "'8.67%"
$ws.Range("E26").Value = "'4.06%"
$ws.Range("D27").Value = "'0.0003395"
$ws.Range("D39").Value = "'0.01570"
$ws.Range("E39").Value = "'-5.70%"
$ws.Range("D40").Value = "'0.04379"
$ws.Range("E40").Value = "'-0.93%"
$ws.Range("D41").Value = "'0.007310"
$ws.Range("E41").Value = "'-1.34%"
$ws.Range("D42").Value = "'0.009833"
$ws.Range("D43").Value = "'0.1312"
$ws.Range("E43").Value = "'-0.96%"
$ws.Range("D44").Value = "'0.002023"
$ws.Range("E44").Value = "'-1.51%"
$ws.Range("D45").Value = "'0.009493"
$ws.Range("E45").Value = "'-14.28%"
$ws.Range("D46").Value = "'0.00005759"
$ws.Range("E46").Value = "'-3.60%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.17%"
$ws.Range("E48").Value = "'6.59%"
$ws.Range("D49").Value = "'0.002900"
$ws.Range("E49").Value = "'-3.51%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.17%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.17%"
